$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is one full data row (columns A through T)
# for the Rarres2 -> Gpr1 ligand-receptor pairs across all 4x4 sending/target cluster combinations.

$row1 = @("ECs", "Rarres2", "Gpr1", "ECs", 3, 1, 1.458525, 4.375575, 0.02525273220298681, 0.02525273220298681, 3, 1, 3.175636, 9.526908, 0.4701396042380186, 0.4701396042380186, 4.6317444969, 41.6857004721, 0.01187230952384089, 0.01187230952384089)
$row2 = @("ECs", "Rarres2", "Gpr1", "FAPs", 3, 1, 1.458525, 4.375575, 0.02525273220298681, 0.02525273220298681, 3, 1, 2.008641666666667, 6.025925, 0.2973709827646055, 0.2973709827646055, 2.929654086875, 26.366886781875, 0.007509429792693589, 0.00750942979269359)
$row3 = @("ECs", "Rarres2", "Gpr1", "MuSCs", 3, 1, 1.458525, 4.375575, 0.02525273220298681, 0.02525273220298681, 3, 1, 1.119552333333333, 3.358657, 0.1657450321501216, 0.1657450321501216, 1.632895066975, 14.696055602775, 0.004185514910862461, 0.004185514910862461)
$row4 = @("ECs", "Rarres2", "Gpr1", "Resolving-Mac", 3, 1, 1.458525, 4.375575, 0.02525273220298681, 0.02525273220298681, 2, 0.6666666666666666, 0.450836, 1.352508, 0.06674438084725434, 0.06674438084725433, 0.6575555769, 5.918000192099999, 0.001685477975589876, 0.001685477975589876)
$row5 = @("FAPs", "Rarres2", "Gpr1", "ECs", 3, 1, 48.53546666666667, 145.6064, 0.84033742450786, 0.8403374245078601, 3, 1, 3.175636, 9.526908, 0.4701396042380186, 0.4701396042380186, 154.1309752234667, 1387.1787770112, 0.3950759041845212, 0.3950759041845212)
$row6 = @("FAPs", "Rarres2", "Gpr1", "FAPs", 3, 1, 48.53546666666667, 145.6064, 0.84033742450786, 0.8403374245078601, 3, 1, 2.008641666666667, 6.025925, 0.2973709827646055, 0.2973709827646055, 97.49036065777779, 877.41324592, 0.2498919657797798, 0.2498919657797798)
$row7 = @("FAPs", "Rarres2", "Gpr1", "MuSCs", 3, 1, 48.53546666666667, 145.6064, 0.84033742450786, 0.8403374245078601, 3, 1, 1.119552333333333, 3.358657, 0.1657450321501216, 0.1657450321501216, 54.33799495608889, 489.0419546048, 0.1392817534420057, 0.1392817534420057)
$row8 = @("FAPs", "Rarres2", "Gpr1", "Resolving-Mac", 3, 1, 48.53546666666667, 145.6064, 0.84033742450786, 0.8403374245078601, 2, 0.6666666666666666, 0.450836, 1.352508, 0.06674438084725434, 0.06674438084725433, 21.88153565013334, 196.9338208512, 0.05608780110155345, 0.05608780110155345)
$row9 = @("MuSCs", "Rarres2", "Gpr1", "ECs", 3, 1, 7.712822, 23.138466, 0.1335389029981009, 0.133538902998101, 3, 1, 3.175636, 9.526908, 0.4701396042380186, 0.4701396042380186, 24.493115204792, 220.438036843128, 0.06278192700590633, 0.06278192700590633)
$row10 = @("MuSCs", "Rarres2", "Gpr1", "FAPs", 3, 1, 7.712822, 23.138466, 0.1335389029981009, 0.133538902998101, 3, 1, 2.008641666666667, 6.025925, 0.2973709827646055, 0.2973709827646055, 15.49229563678333, 139.43066073105, 0.03971059482185259, 0.0397105948218526)
$row11 = @("MuSCs", "Rarres2", "Gpr1", "MuSCs", 3, 1, 7.712822, 23.138466, 0.1335389029981009, 0.133538902998101, 3, 1, 1.119552333333333, 3.358657, 0.1657450321501216, 0.1657450321501216, 8.634907866684665, 77.71417080016201, 0.02213340977071221, 0.02213340977071221)
$row12 = @("MuSCs", "Rarres2", "Gpr1", "Resolving-Mac", 3, 1, 7.712822, 23.138466, 0.1335389029981009, 0.133538902998101, 2, 0.6666666666666666, 0.450836, 1.352508, 0.06674438084725434, 0.06674438084725433, 3.477217819192, 31.294960372728, 0.008912971399629802, 0.008912971399629802)
$row13 = @("Resolving-Mac", "Rarres2", "Gpr1", "ECs", 1, 0.3333333333333333, 0.05030299999999999, 0.150909, 0.0008709402910521559, 0.000870940291052156, 3, 1, 3.175636, 9.526908, 0.4701396042380186, 0.4701396042380186, 0.159744017708, 1.437696159372, 0.0004094635237502053, 0.0004094635237502053)
$row14 = @("Resolving-Mac", "Rarres2", "Gpr1", "FAPs", 1, 0.3333333333333333, 0.05030299999999999, 0.150909, 0.0008709402910521559, 0.000870940291052156, 3, 1, 2.008641666666667, 6.025925, 0.2973709827646055, 0.2973709827646055, 0.1010407017583333, 0.9093663158249999, 0.0002589923702794711, 0.0002589923702794712)
$row15 = @("Resolving-Mac", "Rarres2", "Gpr1", "MuSCs", 1, 0.3333333333333333, 0.05030299999999999, 0.150909, 0.0008709402910521559, 0.000870940291052156, 3, 1, 1.119552333333333, 3.358657, 0.1657450321501216, 0.1657450321501216, 0.05631684102366666, 0.506851569213, 0.0001443540265412758, 0.0001443540265412758)
$row16 = @("Resolving-Mac", "Rarres2", "Gpr1", "Resolving-Mac", 1, 0.3333333333333333, 0.05030299999999999, 0.150909, 0.0008709402910521559, 0.000870940291052156, 2, 0.6666666666666666, 0.450836, 1.352508, 0.06674438084725434, 0.06674438084725433, 0.022678403308, 0.204105629772, 0.00005813037048120363, 0.00005813037048120363)

$dataRows = @($row1,$row2,$row3,$row4,$row5,$row6,$row7,$row8,$row9,$row10,$row11,$row12,$row13,$row14,$row15,$row16)

$r = 2
foreach ($row in $dataRows) {
  $c = 1
  foreach ($val in $row) {
    $ws.Cells.Item($r, $c).Value2 = $val
    $c++
  }
  $r++
}

Write-Host "Updated $($dataRows.Count) data rows (rows 2 through $($r-1))."
